$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B1").Value = "Localizacion"
$ws.Range("D1").Value = "Identificador"
$ws.Range("E1").Value = "Tipo"

# Row 2 (Juan)
$ws.Range("B2").Value = "C/ Federico García Lorca 2"
$ws.Range("D2").Value = "123a"
$ws.Range("E2").Value = "1,entidad"

# Row 3 (Luis)
$ws.Range("B3").Value = "C/ Real Oviedo 2"
$ws.Range("D3").Value = "77b"
$ws.Range("E3").Value = "2,ciudadano"

# Row 4 (Ana)
$ws.Range("B4").Value = "Av. De la Constitución 8"
$ws.Range("D4").Value = "88c"
$ws.Range("E4").Value = "3,sensor"

# D4 previously used a different style (date + bold-font variant) - align it with
# D2/D3's plain date-format style now that all three hold text values
$ws.Range("D4").NumberFormat = "m/d/yy"

# Clear now-unused columns F and G (keep G4's leftover style, matching original artifact)
$ws.Range("F1:G4").ClearContents()

# Update selection to match target (E2 instead of G1)
$ws.Range("E2").Select() | Out-Null
